$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.679.66'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '1.588.29'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').Value = "'207.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('E6').Value = '  -3.29%  '
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('D8').Value = "'22.23"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.50%  '
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = '1.814.28'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').Value = '1.592.96'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('D14').Value = "'3.85"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('D15').Value = "'0.529"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.79%  '
$ws.Range('D16').Value = '27.656.36'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').Value = "'63.42"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.38%  '
$ws.Range('D18').Value = "'219.19"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.12%  '
$ws.Range('D19').Value = '0.0₃0697'
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('D20').Value = "'7.31"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.91%  '
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('E22').Value = '  -4.98%  '
$ws.Range('D23').Value = "'9.58"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.60%  '
$ws.Range('E24').Value = '  -3.89%  '
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('D26').Value = "'6.83"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').Value = "'15.10"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('E29').Value = '  -4.68%  '
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('E32').Value = '  -5.26%  '
$ws.Range('D33').Value = '1.370.12'
$ws.Range('E33').Value = '  -3.12%  '
$ws.Range('E34').Value = '  -5.75%  '
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('D36').Value = "'0.980"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E39').Value = '  -3.18%  '
$ws.Range('E40').Value = '  -3.43%  '
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').Value = "'0.971"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.03%  '
$ws.Range('D43').Value = "'64.08"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('E44').Value = '  +2.59%  '
$ws.Range('E45').Value = '  -3.72%  '
$ws.Range('E46').Value = '  -4.77%  '
$ws.Range('D47').Value = '1.725.52'
$ws.Range('E47').Value = '  -2.51%  '
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').Value = '0.0₆0100'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('E50').Value = '  -4.30%  '
$ws.Range('D51').Value = "'0.0494"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.69%  '
